# This script updates the win-probability transition matrix on the active sheet
# to reflect newly simulated games (updated row-normalized distributions).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Af0)
$ws.Range("B2").Value = 0.2244094488188976
$ws.Range("C2").Value = 0.4881889763779528
$ws.Range("J2").Value = 0.01968503937007874
$ws.Range("P2").Value = 0.1653543307086614
$ws.Range("S2").Value = 0.1023622047244094

# Row 3 (Af1)
$ws.Range("B3").Value = 0.0078125
$ws.Range("C3").Value = 0.0078125
$ws.Range("J3").Value = 0.046875
$ws.Range("P3").Value = 0.65625
$ws.Range("S3").Value = 0.28125

# Row 4 (Af2)
$ws.Range("P4").Value = 0.6875
$ws.Range("S4").Value = 0.3125

# Row 6 (Ai0)
$ws.Range("B6").Value = 0.048
$ws.Range("D6").Value = 0.02
$ws.Range("F6").Value = 0.08
$ws.Range("J6").Value = 0.228
$ws.Range("O6").Value = 0.028
$ws.Range("Q6").Value = 0.14
$ws.Range("R6").Value = 0.048
$ws.Range("S6").Value = 0.408

# Row 7 (Ai1)
$ws.Range("B7").Value = 0.1048034934497817
$ws.Range("D7").Value = 0.01310043668122271
$ws.Range("F7").Value = 0.06550218340611354
$ws.Range("J7").Value = 0.1572052401746725
$ws.Range("O7").Value = 0.03493449781659388
$ws.Range("Q7").Value = 0.1179039301310044
$ws.Range("R7").Value = 0.05676855895196507
$ws.Range("S7").Value = 0.4497816593886463

# Row 8 (Ai2)
$ws.Range("B8").Value = 0.08583690987124463
$ws.Range("D8").Value = 0.006437768240343348
$ws.Range("E8").Value = 0.002145922746781116
$ws.Range("F8").Value = 0.07939914163090128
$ws.Range("J8").Value = 0.08369098712446352
$ws.Range("O8").Value = 0.02575107296137339
$ws.Range("Q8").Value = 0.1287553648068669
$ws.Range("R8").Value = 0.07510729613733906
$ws.Range("S8").Value = 0.5128755364806867

# Row 9 (Ai3)
$ws.Range("B9").Value = 0.05524861878453038
$ws.Range("D9").Value = 0.01104972375690608
$ws.Range("F9").Value = 0.08839779005524862
$ws.Range("J9").Value = 0.1270718232044199
$ws.Range("O9").Value = 0.02762430939226519
$ws.Range("Q9").Value = 0.143646408839779
$ws.Range("R9").Value = 0.04419889502762431
$ws.Range("S9").Value = 0.5027624309392266

# Row 10 (Ar0)
$ws.Range("B10").Value = 0.1098901098901099
$ws.Range("D10").Value = 0.01898101898101898
$ws.Range("F10").Value = 0.0969030969030969
$ws.Range("J10").Value = 0.1348651348651349
$ws.Range("O10").Value = 0.02097902097902098
$ws.Range("Q10").Value = 0.1868131868131868
$ws.Range("R10").Value = 0.03296703296703297
$ws.Range("S10").Value = 0.3986013986013986

# Row 11 (Bf0)
$ws.Range("G11").Value = 0.1480446927374302
$ws.Range("J11").Value = 0.06983240223463687
$ws.Range("K11").Value = 0.1983240223463687
$ws.Range("L11").Value = 0.5418994413407822
$ws.Range("S11").Value = 0.04189944134078212

# Row 12 (Bf1)
$ws.Range("G12").Value = 0.7358490566037735
$ws.Range("J12").Value = 0.1839622641509434
$ws.Range("K12").Value = 0.01415094339622642
$ws.Range("L12").Value = 0.02358490566037736
$ws.Range("S12").Value = 0.04245283018867924

# Row 13 (Bf2)
$ws.Range("G13").Value = 0.6739130434782609
$ws.Range("J13").Value = 0.1956521739130435
$ws.Range("S13").Value = 0.1304347826086956

# Row 15 (Bi0)
$ws.Range("F15").Value = 0.03076923076923077
$ws.Range("H15").Value = 0.1641025641025641
$ws.Range("I15").Value = 0.05641025641025641
$ws.Range("J15").Value = 0.2871794871794872
$ws.Range("K15").Value = 0.07692307692307693
$ws.Range("N15").Value = 0.005128205128205128
$ws.Range("O15").Value = 0.07179487179487179
$ws.Range("S15").Value = 0.3076923076923077

# Row 16 (Bi1)
$ws.Range("F16").Value = 0.03472222222222222
$ws.Range("H16").Value = 0.2152777777777778
$ws.Range("I16").Value = 0.05555555555555555
$ws.Range("J16").Value = 0.3611111111111111
$ws.Range("K16").Value = 0.1180555555555556
$ws.Range("M16").Value = 0.01388888888888889
$ws.Range("O16").Value = 0.04861111111111111
$ws.Range("S16").Value = 0.1527777777777778

# Row 17 (Bi2)
$ws.Range("F17").Value = 0.01176470588235294
$ws.Range("H17").Value = 0.2117647058823529
$ws.Range("I17").Value = 0.1205882352941176
$ws.Range("J17").Value = 0.3441176470588235
$ws.Range("K17").Value = 0.1058823529411765
$ws.Range("M17").Value = 0.01470588235294118
$ws.Range("N17").Value = 0.002941176470588235
$ws.Range("O17").Value = 0.05
$ws.Range("S17").Value = 0.1382352941176471

# Row 18 (Bi3)
$ws.Range("F18").Value = 0.02
$ws.Range("H18").Value = 0.23
$ws.Range("I18").Value = 0.08
$ws.Range("J18").Value = 0.28
$ws.Range("K18").Value = 0.17
$ws.Range("M18").Value = 0.03
$ws.Range("O18").Value = 0.04
$ws.Range("S18").Value = 0.15

# Row 19 (Br0)
$ws.Range("F19").Value = 0.02017937219730942
$ws.Range("H19").Value = 0.2272047832585949
$ws.Range("I19").Value = 0.08445440956651719
$ws.Range("J19").Value = 0.2989536621823617
$ws.Range("K19").Value = 0.1390134529147982
$ws.Range("M19").Value = 0.02765321375186846
$ws.Range("O19").Value = 0.05680119581464873
$ws.Range("S19").Value = 0.1457399103139013
